$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19486
$ws.Range("B3").Value = 14660
$ws.Range("B4").Value = 1826
$ws.Range("B5").Value = 18159
